$d = $word.ActiveDocument

# ------------------------------------------------------------------
# helper text blocks
# ------------------------------------------------------------------
$objetivosPt = "Apresentar os fundamentos da Contabilidade Geral, de modo que os alunos possam compreender os conceitos e princípios básicos de Contabilidade. Capacitar na produção e análise de documentos e de registros de transações contábeis. Compreender os sistemas de custeio e sua aplicação. Capacitar a projetar e implantar sistemas de custeio. Capacitar a analisar custos em relação a estratégia de mercado e de produção de uma organização."
$objetivosEn = "Present the fundamentals of General Accounting, so that students can understand the basic concepts and principles of Accounting. Train in the production and analysis of documents and records of accounting transactions. Understand costing systems and their application. Empower to design and implement costing systems. Empower to analyze costs in relation to an organization's market and production strategy."

$programaResumidoEn = "1 - Fundamentals of Accounting. 2 - Costing Systems. 3 - Costing Systems Project. 4 - Cost Analysis."

$programaPt = "1. FUNDAMENTOS DE CONTABILIDADEConceito de contabilidade. Contabilidade e Engenharia de Produção. Balanço: ativo, passivo, patrimônio líquido. Procedimentos contábeis básicos. Variações da situação líquida; despesa, receita. Regimes de competência e caixa. Receitas e despesas diferidas. Fatos Contábeis. Operações com mercadorias. Inventário. Ativo Imobilizado e Amortização. Demonstrativo de origens e aplicações. Análise de Balanço e de resultados.2. SISTEMAS DE CUSTEIOCustos diretos e indiretos, fixos e variáveis. Acumulação de custos, classificação, fatores de custo. Custeio por Absorção. Produção por ordem, contínua, conjunta. Custeio Direto: margem de contribuição. Custeio ABC.3. PROJETO DE SISTEMAS DE CUSTEIOProdutos e Departamentos. Sistemas de produção e sistemas de custeio. O problema da inflação.4. ANÁLISE DE CUSTOSCusto-Volume-Lucro. Contribuição marginal. Análise de variações. Equação de produtividade global. Alavancagem operacional. TIR e lucratividade."
$programaEn = "1. ACCOUNTING FUNDAMENTALSAccounting concept. Production Accounting and Engineering. Balance sheet: assets, liabilities, equity. Basic accounting procedures. Changes in equity; expense, revenue. Accrual and cash regimes. Deferred income and expenses. Accounting Facts. Goods operations. Inventory. Property, Plant and Equipment and Amortization. Statement of origins and applications. Balance sheet and results analysis.2. COSTING SYSTEMSDirect and indirect costs, fixed and variable. Accumulation of costs, classification, cost factors. Absorption costing. Orderly, continuous, joint production. Direct Costing: contribution margin. ABC costing.3. COSTING SYSTEMS DESIGNProducts and Departments. Production and costing systems. The inflation problem.4. COST ANALYSISCost-Volume-Profit. Marginal contribution. Analysis of variations. Global productivity equation. Operational leverage. IRR and profitability."

$bibNew = "IUDICIBUS, S.; MARION, J. C. Curso de Contabilidade para não Contadores. 8 ed. São Paulo: Atlas, 2018.MARTINS, E. Contabilidade de Custos. 11 ed. São Paulo: Atlas, 2018.SANTOS, J. J. Manual de Contabilidade e Análise de Custos. 7 ed. São Paulo: Atlas, 2017.SILVA, R. N. S.; LINS, L. S. Gestão de Custos - Contabilidade, Controle e Análise. 4 ed. São Paulo: Atlas, 2017.DUTRA, R. G. Custos: Uma Abordagem Prática. 8 ed. São Paulo: Atlas, 2017.MEGLIORINI, E.; BUENO, A. S. Contabilidade para cursos de Engenharia. São Paulo: Atlas, 2014.MORANTE, A. S. Análise das Demonstrações Financeiras. 2 ed. São Paulo: Atlas, 2009.Bruni, Adriano L., Fama, Rubens. Gestão de custos e formação de preços. 5. Ed., São Paulo: Atlas, 2008.Hansen, Don R., Mowen, Maryanne M., Gestão de custos – contabilidade e controle. São Paulo: Pioneira Thomson, 2001Horngren, Charles T., Datar, Srikant M., Foster, George. Contabilidade de custos (vol. 1 e 2). 11. Ed. São Paulo: Pearson, 2004.Maher, Michael. Contabilidade de custos – criando valor para a administração. 5. ed., São Paulo: Atlas, 2001.Martins, Eliseu. Contabilidade de Custos (livro de exercícios). 9. Ed., São Paulo: Atlas, 2006.Martins, Eliseu. Contabilidade de Custos (livro texto). 9. Ed., São Paulo: Atlas, 2003."

# ------------------------------------------------------------------
# 1. Activation date update
#    (runs around "Departamento: ..." get coalesced into the edited
#    run on save if their formatting matches, so temporarily blank
#    the next run, perform the edit, then re-insert its text via an
#    insertion point -- that creates a fresh run instead of
#    re-joining the previous one)
# ------------------------------------------------------------------
$deptRng = $d.Content
$deptRng.Find.Execute("Departamento: Engenharia Química", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$deptRng.Text = ""

$d.Content.Find.Execute(
    "Ativação: 01/01/2012", $true, $false, $false, $false, $false,
    $true, 1, $false, "Ativação: 01/01/2021", 2) | Out-Null

$ativRng = $d.Content
$ativRng.Find.Execute("Ativação: 01/01/2021", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ativPara = $ativRng.Paragraphs(1)
$deptInsertionPoint = $d.Range($ativPara.Range.End - 1, $ativPara.Range.End - 1)
$deptInsertionPoint.InsertAfter("Departamento: Engenharia Química")

# ------------------------------------------------------------------
# 2. Objetivos: replace paragraph text, then add italic EN paragraph
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Fornecer ao aluno o conhecimento da contabilidade de custos de forma a capacitá-lo a classificar, apurar contabilizar e interpretar informações de custos.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$objPara = $rng.Paragraphs(1)
$objPara.Range.Text = $objetivosPt

$objPara.Range.InsertParagraphAfter()
$objEnPara = $objPara.Next()
$objEnStart = $objEnPara.Range.Start
$objEnPara.Range.InsertAfter($objetivosEn)
$objEnRange = $d.Range($objEnStart, $objEnStart + $objetivosEn.Length)
$objEnRange.Font.Italic = $true

# ------------------------------------------------------------------
# 3. Docente responsável update
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "5840671 - Francisco José Moreira Chaves", $true, $false, $false, $false, $false,
    $true, 1, $false, "11079086 - Herlandí de Souza Andrade", 2) | Out-Null

# ------------------------------------------------------------------
# 4. Programa resumido: add italic EN paragraph right after it
# ------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("1  Fundamentos de Contabilidade. 2 - Sistemas de Custeio. 3 - Projeto de Sistemas de Custeio. 4 - Análise de Custos.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$resumoPara = $rng2.Paragraphs(1)
$resumoPara.Range.InsertParagraphAfter()
$resumoEnPara = $resumoPara.Next()
$resumoEnStart = $resumoEnPara.Range.Start
$resumoEnPara.Range.InsertAfter($programaResumidoEn)
$resumoEnRange = $d.Range($resumoEnStart, $resumoEnStart + $programaResumidoEn.Length)
$resumoEnRange.Font.Italic = $true

# ------------------------------------------------------------------
# 5. Programa: merge the body text into one run, then add italic EN
# ------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("1. FUNDAMENTOS DE CONTABILIDADE", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$progPara = $rng3.Paragraphs(1)
$progPara.Range.Text = $programaPt

$progPara.Range.InsertParagraphAfter()
$progEnPara = $progPara.Next()
$progEnStart = $progEnPara.Range.Start
$progEnPara.Range.InsertAfter($programaEn)
$progEnRange = $d.Range($progEnStart, $progEnStart + $programaEn.Length)
$progEnRange.Font.Italic = $true

# ------------------------------------------------------------------
# 6. Avaliação: Método / Critério / Norma de recuperação updates
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Aulas expositivas com aplicação de exercícios, trabalhos e provas.", $true, $false, $false, $false, $false,
    $true, 1, $false, "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras.", 2) | Out-Null

$d.Content.Find.Execute(
    "MF=(0,8xP+0,2xET) , onde: MF = Média Final da primeira avaliação; P = média ponderada das provas; ET = média ponderada dos exercícios e trabalhos.", $true, $false, $false, $false, $false,
    $true, 1, $false, "Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas.", 2) | Out-Null

$d.Content.Find.Execute(
    "NF=(MF+REC)/2, onde: NF = Média Final da segunda avaliação e REC = nota obtida na prova do período de recuperação.", $true, $false, $false, $false, $false,
    $true, 1, $false, "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.", 2) | Out-Null

# ------------------------------------------------------------------
# 7. Bibliografia: merge into a single run, new refs prepended
# ------------------------------------------------------------------
$rng4 = $d.Content
$rng4.Find.Execute("Bruni, Adriano L., Fama, Rubens.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bibPara = $rng4.Paragraphs(1)
$bibPara.Range.Text = $bibNew
